# Changes to BOM for fab:
#  - Add fabricator notes: "DO NOT PLACE. No Solder paste" (G12) and
#    "4 pin-header" (G17), matching the comment style already used in
#    column G (copy format from G15, the existing "No Solder paste..." note).
#  - Update the sheet view: zoom to 160%, clear the stale topLeftCell,
#    and leave the selection on G17 (the cell just edited).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024-Badge-BOM-REV1")

# --- Cell content updates -------------------------------------------------

$commentFormat = $ws.Range("G15")

$g12 = $ws.Range("G12")
$g12.Value = "DO NOT PLACE. No Solder paste"
$commentFormat.Copy()
$g12.PasteSpecial(-4122) # xlPasteFormats

$g17 = $ws.Range("G17")
$g17.Value = "4 pin-header"
$commentFormat.Copy()
$g17.PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = 0

# --- Sheet view updates ----------------------------------------------------

$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 160
$g17.Select()
